$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value  = 12.741
$ws.Range("D3").Value  = -7.592000000000001
$ws.Range("B4").Value  = 5.590999999999999
$ws.Range("C4").Value  = -13.36
$ws.Range("D4").Value  = -8.007999999999999
$ws.Range("C5").Value  = -12.927
$ws.Range("B6").Value  = 6.741
$ws.Range("B7").Value  = 6.612
$ws.Range("C8").Value  = -12.883
$ws.Range("D9").Value  = -8.015000000000001
$ws.Range("D11").Value = -7.653
$ws.Range("D14").Value = -8.028
$ws.Range("B16").Value = 6.271
$ws.Range("C16").Value = -12.673
$ws.Range("D18").Value = -7.637
$ws.Range("B20").Value = 6.064000000000001
$ws.Range("E20").Value = 12.761
$ws.Range("C22").Value = -12.504
$ws.Range("D25").Value = -7.653
